# Update the heading date/day paragraph.
$d = $word.ActiveDocument
$d.Paragraphs.Item(1).Range.Text = "2023-11-16 Thursday"

# Replace the arithmetic problem in every cell of the practice table.
# Each row of $values holds the 5 new cell strings (left to right) for
# that table row, matching the 20x5 grid in the document.
$tbl = $d.Tables.Item(1)

$values = @(
    @("49+19=", "20-4=", "71-36=", "42-3=", "6+66="),
    @("8+74=", "94-47=", "6+17=", "15+57=", "68+24="),
    @("80-27=", "16+46=", "40-2=", "70-27=", "62-7="),
    @("23+29=", "56-29=", "94-48=", "69+2=", "40-24="),
    @("69+3=", "47+8=", "79+18=", "27+26=", "95-39="),
    @("81-36=", "5+46=", "51-18=", "81-32=", "94-26="),
    @("80-29=", "59+23=", "73-15=", "77+7=", "35-19="),
    @("94-49=", "6+8=", "27+17=", "74-16=", "16-7="),
    @("44+27=", "17+76=", "53-27=", "80-73=", "14+28="),
    @("69+22=", "54-25=", "53-45=", "25+18=", "71-8="),
    @("9+12=", "66-17=", "20-14=", "51-27=", "92-28="),
    @("77+19=", "91-9=", "46+19=", "52-44=", "91-53="),
    @("41-37=", "25+36=", "30-18=", "60-14=", "54-6="),
    @("32-27=", "49+22=", "79+15=", "78-69=", "97-69="),
    @("77-58=", "52-6=", "67+27=", "39+45=", "28+48="),
    @("94-89=", "27+6=", "48+45=", "48-9=", "40-35="),
    @("29+42=", "63-16=", "48-19=", "59+24=", "97-49="),
    @("81-65=", "18+77=", "78+18=", "87+6=", "79+3="),
    @("85-76=", "24+37=", "28+33=", "18+4=", "38+9="),
    @("24+39=", "31-26=", "8+7=", "26+27=", "62-33=")
)

for ($r = 1; $r -le $values.Count; $r++) {
    $rowValues = $values[$r - 1]
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $tbl.Cell($r, $c).Range.Text = $rowValues[$c - 1]
    }
}
